# Weekly data refresh for "Hortaliza, Vega Modelo de Temuco - Coliflor".
#
# The historic series occupies rows 638-745 (18 columns, A:R). Two new
# observations are prepended at rows 638-639; every existing row from the
# old 638 downward shifts two rows further down (638->640, 639->641, ...,
# 745->747), so the sheet grows from R745 to R747.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 638
$lastRow  = 745
$numCols  = 18   # A..R
$shift    = 2

# Read the whole existing block (rows 638..745, columns A..R) in one shot.
# NOTE: Range.Value2 returns a 1-based 2D array ($old[1,1] .. $old[n,m]).
$srcRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, $numCols))
$old = $srcRange.Value2

$oldRowCount = $lastRow - $firstRow + 1
$newRowCount = $oldRowCount + $shift

# Build the shifted block in memory (0-based .NET array, since Range.Value2
# assignment just needs a correctly-sized rectangular array):
#   newBlock[r, c] (0-based) == sheet row (firstRow + r), column (c+1)
# newBlock row (i + shift) = old row i, for i = 1..oldRowCount (old is 1-based).
$newBlock = New-Object 'object[,]' $newRowCount, $numCols
for ($i = 1; $i -le $oldRowCount; $i++) {
    for ($j = 1; $j -le $numCols; $j++) {
        $newBlock[($i + $shift) - 1, $j - 1] = $old[$i, $j]
    }
}

# The first two rows of the block (new sheet rows 638 and 639) are brand
# new observations; seed them from the template of the old first/second row
# (same Mercado/Region/Categoria/etc.), then overwrite the fields that
# actually carry new data.
for ($j = 1; $j -le $numCols; $j++) {
    $newBlock[0, $j - 1] = $old[1, $j]
    $newBlock[1, $j - 1] = $old[2, $j]
}

# Row 638 (new) - columns: D=4 J=10 K=11 L=12 M=13 O=15 P=16
$newBlock[0, 3]  = 45218
$newBlock[0, 9]  = 680
$newBlock[0, 10] = 1300
$newBlock[0, 11] = 1300
$newBlock[0, 12] = 1300
$newBlock[0, 14] = "Provincia del Elquí"
$newBlock[0, 15] = 1300

# Row 639 (new) - columns: D=4 J=10 K=11 L=12 M=13 O=15 P=16
$newBlock[1, 3]  = 45218
$newBlock[1, 9]  = 2700
$newBlock[1, 10] = 1200
$newBlock[1, 11] = 1300
$newBlock[1, 12] = 1231
$newBlock[1, 14] = "Región del Maule"
$newBlock[1, 15] = 1231

# Write the whole rebuilt block back (rows 638..747).
$dstLastRow = $firstRow + $newRowCount - 1
$dstRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($dstLastRow, $numCols))
$dstRange.Value2 = $newBlock

# The two brand-new trailing rows (746, 747) sit past the sheet's previous
# extent, so column D (the date column) doesn't automatically inherit the
# "date" number format the rest of the D column uses - copy it over
# explicitly so D746/D747 match D638..D745.
$dateFormat = $ws.Cells.Item($lastRow, 4).NumberFormat
$ws.Cells.Item($lastRow + 1, 4).NumberFormat = $dateFormat
$ws.Cells.Item($lastRow + 2, 4).NumberFormat = $dateFormat
